# Update "想去人数" (number of people wanting to go) values in column F
# on sheet "展览" and sheet "全部类型", per the commit's refreshed data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7896
$ws1.Range("F5").Value = 5761
$ws1.Range("F6").Value = 482
$ws1.Range("F11").Value = 305
$ws1.Range("F12").Value = 64

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7896
$ws4.Range("F5").Value = 5761
$ws4.Range("F6").Value = 482
$ws4.Range("F14").Value = 305
$ws4.Range("F15").Value = 64
